$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Cells.Item(5, 2).Value = 6428338
$ws.Cells.Item(5, 3).Value = 'Denmark Division 1'
$ws.Cells.Item(5, 4).Value = 45081.33333333334
$ws.Cells.Item(5, 5).Value = 'Sonderjyske'
$ws.Cells.Item(5, 6).Value = 'FC Helsingor'
$ws.Cells.Item(5, 7).Value = 3
$ws.Cells.Item(5, 8).Value = 2
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 'H'
$ws.Cells.Item(5, 12).Value = 1.5
$ws.Cells.Item(5, 13).Value = 4.5
$ws.Cells.Item(5, 14).Value = 5
$ws.Cells.Item(5, 15).Value = 1.363
$ws.Cells.Item(5, 16).Value = 5.25
$ws.Cells.Item(5, 17).Value = 6.5
$ws.Cells.Item(5, 18).Value = -1.5
$ws.Cells.Item(5, 19).Value = 1.95
$ws.Cells.Item(5, 20).Value = 1.9
$ws.Cells.Item(5, 21).Value = 3.75
$ws.Cells.Item(5, 22).Value = 1.9
$ws.Cells.Item(5, 23).Value = 1.95
$ws.Cells.Item(5, 24).Value = 0.363
$ws.Cells.Item(5, 25).Value = -1
$ws.Cells.Item(5, 26).Value = -1
$ws.Cells.Item(5, 27).Value = -1
$ws.Cells.Item(5, 28).Value = 0.8999999999999999
$ws.Cells.Item(5, 29).Value = 0.8999999999999999
$ws.Cells.Item(5, 30).Value = -1

# Row 6
$ws.Cells.Item(6, 2).Value = 6428339
$ws.Cells.Item(6, 3).Value = 'Denmark Division 1'
$ws.Cells.Item(6, 4).Value = 45081.33333333334
$ws.Cells.Item(6, 5).Value = 'Vejle'
$ws.Cells.Item(6, 6).Value = 'Vendsyssel FF'
$ws.Cells.Item(6, 7).Value = 4
$ws.Cells.Item(6, 8).Value = 3
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 'H'
$ws.Cells.Item(6, 12).Value = 1.55
$ws.Cells.Item(6, 13).Value = 4
$ws.Cells.Item(6, 14).Value = 5
$ws.Cells.Item(6, 15).Value = 2.2
$ws.Cells.Item(6, 16).Value = 3.8
$ws.Cells.Item(6, 17).Value = 3
$ws.Cells.Item(6, 18).Value = -0.25
$ws.Cells.Item(6, 19).Value = 1.875
$ws.Cells.Item(6, 20).Value = 1.975
$ws.Cells.Item(6, 21).Value = 3
$ws.Cells.Item(6, 22).Value = 2.05
$ws.Cells.Item(6, 23).Value = 1.8
$ws.Cells.Item(6, 24).Value = 1.2
$ws.Cells.Item(6, 25).Value = -1
$ws.Cells.Item(6, 26).Value = -1
$ws.Cells.Item(6, 27).Value = 0.875
$ws.Cells.Item(6, 28).Value = -1
$ws.Cells.Item(6, 29).Value = 1.05
$ws.Cells.Item(6, 30).Value = -1

# Row 7
$ws.Cells.Item(7, 2).Value = 6428337
$ws.Cells.Item(7, 3).Value = 'Denmark Division 1'
$ws.Cells.Item(7, 4).Value = 45081.33333333334
$ws.Cells.Item(7, 5).Value = 'Hvidovre IF'
$ws.Cells.Item(7, 6).Value = 'Naestved'
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(7, 8).Value = 3
$ws.Cells.Item(7, 9).Value = 2
$ws.Cells.Item(7, 10).Value = 2
$ws.Cells.Item(7, 11).Value = 'A'
$ws.Cells.Item(7, 12).Value = 1.7
$ws.Cells.Item(7, 13).Value = 4
$ws.Cells.Item(7, 14).Value = 4.2
$ws.Cells.Item(7, 15).Value = 1.571
$ws.Cells.Item(7, 16).Value = 4.75
$ws.Cells.Item(7, 17).Value = 4.75
$ws.Cells.Item(7, 18).Value = -1
$ws.Cells.Item(7, 19).Value = 1.975
$ws.Cells.Item(7, 20).Value = 1.875
$ws.Cells.Item(7, 21).Value = 3.5
$ws.Cells.Item(7, 22).Value = 2
$ws.Cells.Item(7, 23).Value = 1.85
$ws.Cells.Item(7, 24).Value = -1
$ws.Cells.Item(7, 25).Value = -1
$ws.Cells.Item(7, 26).Value = 3.75
$ws.Cells.Item(7, 27).Value = -1
$ws.Cells.Item(7, 28).Value = 0.875
$ws.Cells.Item(7, 29).Value = 1
$ws.Cells.Item(7, 30).Value = -1

# Row 40
$ws.Cells.Item(40, 2).Value = 6798552
$ws.Cells.Item(40, 3).Value = 'Denmark Division 1'
$ws.Cells.Item(40, 4).Value = 45161.58333333334
$ws.Cells.Item(40, 5).Value = 'Vendsyssel FF'
$ws.Cells.Item(40, 6).Value = 'Kolding IF'
$ws.Cells.Item(40, 7).Value = 2
$ws.Cells.Item(40, 8).Value = 1
$ws.Cells.Item(40, 9).Value = 1
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 'H'
$ws.Cells.Item(40, 12).Value = 2.2
$ws.Cells.Item(40, 13).Value = 3.5
$ws.Cells.Item(40, 14).Value = 2.8
$ws.Cells.Item(40, 15).Value = 2.7
$ws.Cells.Item(40, 16).Value = 3.5
$ws.Cells.Item(40, 17).Value = 2.5
$ws.Cells.Item(40, 18).Value = 0
$ws.Cells.Item(40, 19).Value = 2
$ws.Cells.Item(40, 20).Value = 1.85
$ws.Cells.Item(40, 21).Value = 2.5
$ws.Cells.Item(40, 22).Value = 1.875
$ws.Cells.Item(40, 23).Value = 1.975
$ws.Cells.Item(40, 24).Value = 1.7
$ws.Cells.Item(40, 25).Value = -1
$ws.Cells.Item(40, 26).Value = -1
$ws.Cells.Item(40, 27).Value = 1
$ws.Cells.Item(40, 28).Value = -1
$ws.Cells.Item(40, 29).Value = 0.875
$ws.Cells.Item(40, 30).Value = -1

# Row 42
$ws.Cells.Item(42, 2).Value = 6799263
$ws.Cells.Item(42, 3).Value = 'Denmark Division 1'
$ws.Cells.Item(42, 4).Value = 45161.58333333334
$ws.Cells.Item(42, 5).Value = 'FC Helsingor'
$ws.Cells.Item(42, 6).Value = 'Hillerd'
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = 6
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 2
$ws.Cells.Item(42, 11).Value = 'A'
$ws.Cells.Item(42, 12).Value = 1.75
$ws.Cells.Item(42, 13).Value = 3.8
$ws.Cells.Item(42, 14).Value = 4
$ws.Cells.Item(42, 15).Value = 1.75
$ws.Cells.Item(42, 16).Value = 4
$ws.Cells.Item(42, 17).Value = 4.333
$ws.Cells.Item(42, 18).Value = -0.75
$ws.Cells.Item(42, 19).Value = 2
$ws.Cells.Item(42, 20).Value = 1.85
$ws.Cells.Item(42, 21).Value = 3
$ws.Cells.Item(42, 22).Value = 1.925
$ws.Cells.Item(42, 23).Value = 1.925
$ws.Cells.Item(42, 24).Value = -1
$ws.Cells.Item(42, 25).Value = -1
$ws.Cells.Item(42, 26).Value = 3.333
$ws.Cells.Item(42, 27).Value = -1
$ws.Cells.Item(42, 28).Value = 0.8500000000000001
$ws.Cells.Item(42, 29).Value = 0.925
$ws.Cells.Item(42, 30).Value = -1
